# Daily attendance processing - 2025-12-09 12:47:40
# Reorders "Recorded By" contributor lists, updates rolled-up statistics,
# and flips the PARASITOLOGY SGD/POS session #2 (row 20) from
# "Pending" to "Recorded" now that its attendance has come in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: reorder "Recorded By" list ---
$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System"

# --- Row 3: reorder "Recorded By" list ---
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

# --- Row 4: reorder "Recorded By" list ---
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# --- Row 5: reorder "Recorded By" list ---
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Row 6: reorder "Recorded By" list ---
$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"

# --- Recorded Sessions count (Class Statistics) ---
$ws.Range("L6").Value = 25

# --- Row 7: reorder "Recorded By" list ---
$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg"

# --- Pending Sessions count (Class Statistics) ---
$ws.Range("L8").Value = 2

# --- Coverage % (Class Statistics) ---
# Force text entry (not an auto-converted percentage number), then restore
# the plain "General" look of the surrounding stats cells (s=4) by pasting
# that formatting back on top - the stored value stays the literal text.
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "86.2%"
$ws.Range("L7").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Avg Attendance % (Class Statistics) ---
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "26.7%"
$ws.Range("L7").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 12: reorder "Recorded By" list ---
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# --- Row 15: reorder "Recorded By" list ---
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"

# --- Group Statistics row (row 15, cols O/Q/R/S) ---
$ws.Range("O15").Value = 25
$ws.Range("Q15").Value = 2

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "86.2%"
$ws.Range("L7").Copy()
$ws.Range("R15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "26.7%"
$ws.Range("L7").Copy()
$ws.Range("S15").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 20: PARASITOLOGY SGD/POS session 2 moves from Pending to Recorded ---
# Copy the formatting (fill colour / style) from an already-"Recorded" (green)
# row onto row 20 so its style matches the other recorded rows, without
# disturbing its existing cell values.
$src = $ws.Range("A6:I6")
$dst = $ws.Range("A20:I20")
$src.Copy()
$dst.PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("G20").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("H20").Value = "25/251"
$ws.Range("I20").Value = "Recorded"

# --- Row 25: reorder "Recorded By" list ---
$ws.Range("G25").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"

# --- Row 27: reorder "Recorded By" list ---
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

# --- Row 30: reorder "Recorded By" list ---
$ws.Range("G30").Value = "shorokmohamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
